# Replace the two-row usuario/nombre sample data with the new id/id/letter
# sample rows used for testing the MySQL import (commit: "ok importar xlsx a
# mysql"), add a 4th (empty) column, and re-narrow column C to fit the new,
# much shorter values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new data (row 1 & 2) ---------------------------------------------
$ws.Range("A1").Value = 21
$ws.Range("B1").Value = 21
$ws.Range("C1").Value = "A"

$ws.Range("A2").Value = 22
$ws.Range("B2").Value = 22
$ws.Range("C2").Value = "B"

# --- extend the used range with a blank, but styled, D1/C3 cell -------
# Re-use the same "blank row" style already applied to A3/B3 instead of
# minting a brand new style entry.
$ws.Range("A3").Copy()
$ws.Range("D1").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)

# --- column C is now much narrower since it only holds "A"/"B" --------
$ws.Columns("C").ColumnWidth = 7.15

# --- selection moves to B3 ---------------------------------------------
$ws.Range("B3").Select()
